$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1171.85
$ws.Range("E2").Value = 15.66
$ws.Range("F2").Value = 92.4
$ws.Range("G2").Value = 5.9
$ws.Range("H2").Value = 76.74000000000001
$ws.Range("I2").Value = 1290.911656375089

$ws.Range("D3").Value = 799.17
$ws.Range("E3").Value = 40
$ws.Range("F3").Value = 107.12
$ws.Range("G3").Value = 2.68
$ws.Range("H3").Value = 67.12
$ws.Range("I3").Value = 903.3062832407608

$ws.Range("D4").Value = 592.29
$ws.Range("E4").Value = 41.46
$ws.Range("F4").Value = 81
$ws.Range("G4").Value = 1.95
$ws.Range("H4").Value = 39.54
$ws.Range("I4").Value = 653.6360762714494

$ws.Range("D5").Value = 943.48
$ws.Range("E5").Value = 53.18
$ws.Range("F5").Value = 118.48
$ws.Range("G5").Value = 2.23
$ws.Range("H5").Value = 65.30000000000001
$ws.Range("I5").Value = 1044.792563999131

$ws.Range("D6").Value = 292.19
$ws.Range("E6").Value = 23.46
$ws.Range("F6").Value = 28.32
$ws.Range("G6").Value = 1.21
$ws.Range("H6").Value = 4.859999999999999
$ws.Range("I6").Value = 299.7302612716045

$ws.Range("D7").Value = 127.5
$ws.Range("E7").Value = 46.58
$ws.Range("F7").Value = 14.4
$ws.Range("G7").Value = 0.31
$ws.Range("H7").Value = -32.18

$ws.Range("D8").Value = 820.69
$ws.Range("E8").Value = 16.01
$ws.Range("F8").Value = 62.8
$ws.Range("G8").Value = 3.92
$ws.Range("H8").Value = 46.78999999999999
$ws.Range("I8").Value = 893.2844084153039

$ws.Range("D9").Value = 1938.46
$ws.Range("E9").Value = 59.73
$ws.Range("F9").Value = 205.56
$ws.Range("G9").Value = 3.44
$ws.Range("H9").Value = 145.83
$ws.Range("I9").Value = 2164.714382970801

$ws.Range("D10").Value = 439.45
$ws.Range("E10").Value = 40.33
$ws.Range("F10").Value = 127.48
$ws.Range("G10").Value = 3.16
$ws.Range("H10").Value = 87.15
$ws.Range("I10").Value = 574.6627098395754

$ws.Range("D11").Value = 167.5
$ws.Range("E11").Value = 31.41
$ws.Range("F11").Value = 49.12
$ws.Range("G11").Value = 1.56
$ws.Range("H11").Value = 17.71
$ws.Range("I11").Value = 194.9769603127812

$ws.Range("D12").Value = 684.15
$ws.Range("E12").Value = 23.07
$ws.Range("F12").Value = 68
$ws.Range("G12").Value = 2.95
$ws.Range("H12").Value = 44.93
$ws.Range("I12").Value = 753.8586294101219

$ws.Range("D13").Value = 177.17
$ws.Range("E13").Value = 25.64
$ws.Range("F13").Value = 44.88
$ws.Range("G13").Value = 1.75
$ws.Range("H13").Value = 19.24
$ws.Range("I13").Value = 207.0207462686567

$ws.Range("D14").Value = 100
$ws.Range("E14").Value = 34.29
$ws.Range("F14").Value = 55.8
$ws.Range("G14").Value = 1.63
$ws.Range("H14").Value = 21.51
$ws.Range("I14").Value = 133.3726378502498

$ws.Range("D15").Value = 326.62
$ws.Range("E15").Value = 38.58
$ws.Range("F15").Value = 46.4
$ws.Range("G15").Value = 1.2
$ws.Range("H15").Value = 7.82
$ws.Range("I15").Value = 338.7526837744748
